# Update SwaadSutra_Consolidated_2026-01-21.xlsx
# Adds new order #27 (Renu) as the newest row at the top of the "All Orders"
# log, and a matching new daily-summary row for 2026-01-21 at the top of the
# "Daily Summary" sheet. All existing rows shift down by one.

$wb = $excel.ActiveWorkbook

# ---- Sheet "All Orders": insert new order row at row 2 ----
$ws1 = $wb.Worksheets.Item("All Orders")
$ws1.Rows.Item(2).Insert()

$ws1.Range("A2").Value = 27
$ws1.Range("B2").Value = "2026-01-21 07:49"
$ws1.Range("C2").Value = "Renu"
$ws1.Range("D2").Value = "A-1005 Kakkad la vida"
$ws1.Range("E2").Value = "'8806022013"
$ws1.Range("F2").Value = "Appe Chutney x1, Vermicelli Kheer x1"
$ws1.Range("G2").Value = 110
$ws1.Range("H2").Value = "NEW"
$ws1.Range("I2").Value = "PENDING"
$ws1.Range("J2").Value = "'2026-01-21"
$ws1.Range("K2").Value = "18:30"
$ws1.Range("L2").Value = "Less spicy"
$ws1.Range("M2").Value = "'"
$ws1.Range("N2").Value = "'"

# ---- Sheet "Daily Summary": insert new summary row at row 2 ----
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Rows.Item(2).Insert()

$ws2.Range("A2").Value = "'2026-01-21"
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 110
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 110

Write-Host "Applied SwaadSutra 2026-01-21 update"
